# Update YEAR_START and YEAR_END values on the STANDARD_DEFINITION sheet
# from 1900/2020 to 1700/2040.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("STANDARD_DEFINITION")

$ws.Range("C2:C4").Value = 1700
$ws.Range("D2:D4").Value = 2040

# Make STANDARD_DEFINITION the active sheet and leave the selection on G4,
# matching the author's on-screen state when the file was saved.
$ws.Activate() | Out-Null
$ws.Range("G4").Select() | Out-Null
